$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix company name data entries that were mis-entered (per commit message:
# "Fix data formating from excel and g-spread")
$ws.Range("C5").Value = "Employees Dream Inc."
$ws.Range("C8").Value = "Employees Dream Inc."

# Update the active selection to reflect where the user ended up after editing
$ws.Range("C8").Select()
